$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# Rewrite the test-case rows (2-7) with the updated / reordered content.
# Cells are assigned in strict row-major, left-to-right order so that the
# workbook's shared-string table is rebuilt in the same order as the
# target file (surviving strings keep their original relative order,
# brand-new strings are appended in first-touch order).
# -------------------------------------------------------------------------

# Row 2: ST_SEARCH_05 (new test case)
$ws.Cells.Item(2,1).Value = "ST_SEARCH_05"
$ws.Cells.Item(2,2).Value = "Tìm kiếm không phân biệt hoa thường"
$ws.Cells.Item(2,3).Value = "1. Nhập 'áo hoodie' (chữ thường)`n2. Enter`n3. Check kết quả"
$ws.Cells.Item(2,4).Value = "Keyword: áo hoodie"
$ws.Cells.Item(2,5).Value = "Vẫn tìm thấy sản phẩm (Số lượng > 0)"
$ws.Cells.Item(2,6).Value = "Tìm thấy: 1 sản phẩm."
$ws.Cells.Item(2,7).Value = "PASS"

# Row 3: ST_SEARCH_04 (existing, text tweaked)
$ws.Cells.Item(3,1).Value = "ST_SEARCH_04"
$ws.Cells.Item(3,2).Value = "Tìm kiếm qua Menu Mục lục"
$ws.Cells.Item(3,3).Value = "1. Click Menu 'BỘ SƯU TẬP'`n2. Click mục con 'Áo Hoodie'`n3. Check kết quả"
$ws.Cells.Item(3,4).Value = "Chọn: Áo Hoodie"
$ws.Cells.Item(3,5).Value = "Chuyển trang tìm kiếm & Hiện sản phẩm liên quan (Số lượng > 0)"
$ws.Cells.Item(3,6).Value = "URL hiện tại: http://localhost:8080/ShopDuck/user/search-products?keyword=%C3%81o%20Hoodie | Số SP tìm thấy: 1"
$ws.Cells.Item(3,7).Value = "PASS"

# Row 4: ST_SEARCH_06 (new test case)
$ws.Cells.Item(4,1).Value = "ST_SEARCH_06"
$ws.Cells.Item(4,2).Value = "Tìm kiếm ký tự đặc biệt"
$ws.Cells.Item(4,3).Value = "1. Nhập '@#`$'`n2. Enter`n3. Check không bị lỗi Server (500)"
$ws.Cells.Item(4,4).Value = "Keyword: @#`$"
$ws.Cells.Item(4,5).Value = "Hệ thống xử lý an toàn (Không bị Crash/Lỗi Server)"
$ws.Cells.Item(4,6).Value = "Kiểm tra lỗi Server..."
$ws.Cells.Item(4,7).Value = "PASS"

# Row 5: ST_SEARCH_01 (existing, text tweaked)
$ws.Cells.Item(5,1).Value = "ST_SEARCH_01"
$ws.Cells.Item(5,2).Value = "Tìm kiếm có kết quả (Gõ phím)"
$ws.Cells.Item(5,3).Value = "1. Nhập keyword 'Áo' vào ô search`n2. Enter`n3. Check URL và Số lượng SP"
$ws.Cells.Item(5,4).Value = "Keyword: Áo"
$ws.Cells.Item(5,5).Value = "Hiển thị danh sách sản phẩm liên quan đến 'Áo' (Số lượng > 0)"
$ws.Cells.Item(5,6).Value = "URL: http://localhost:8080/ShopDuck/user/search-products?keyword=%C3%81o | Tìm thấy: 12 sản phẩm."
$ws.Cells.Item(5,7).Value = "PASS"

# Row 6: ST_SEARCH_03 (existing, text tweaked) - newly appended row
$ws.Cells.Item(6,1).Value = "ST_SEARCH_03"
$ws.Cells.Item(6,2).Value = "Bỏ trống từ khóa (Validation)"
$ws.Cells.Item(6,3).Value = "1. Để trống ô search`n2. Nhấn nút Tìm (button)`n3. Check URL không đổi"
$ws.Cells.Item(6,4).Value = "Keyword: (rỗng)"
$ws.Cells.Item(6,5).Value = "Trình duyệt chặn submit, URL không thay đổi"
$ws.Cells.Item(6,6).Value = "URL sau khi click: http://localhost:8080/ShopDuck/user/view-products"
$ws.Cells.Item(6,7).Value = "PASS"

# Row 7: ST_SEARCH_02 (existing, text tweaked) - newly appended row
$ws.Cells.Item(7,1).Value = "ST_SEARCH_02"
$ws.Cells.Item(7,2).Value = "Tìm kiếm không có kết quả"
$ws.Cells.Item(7,3).Value = "1. Nhập từ khóa rác 'xyz123'`n2. Enter`n3. Check thông báo lỗi và danh sách rỗng"
$ws.Cells.Item(7,4).Value = "Keyword: xyz123_khong_co_dau"
$ws.Cells.Item(7,5).Value = "Hiển thị thông báo 'Không tìm thấy' / Danh sách rỗng"
$ws.Cells.Item(7,6).Value = "Thông báo: Không tìm thấy sản phẩm nào phù hợp"
$ws.Cells.Item(7,7).Value = "PASS"

# The G column of the two newly-added rows needs the same "PASS" green/bold
# style used by the rest of the column; copy it across from an existing cell.
$ws.Cells.Item(2,7).Copy()
$ws.Cells.Item(6,7).PasteSpecial(-4122)
$ws.Cells.Item(7,7).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -------------------------------------------------------------------------
# Column width adjustments (B, C, E) to match the regenerated "best fit"
# widths for the new content.
# -------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 34.5
$ws.Columns.Item(3).ColumnWidth = 36.666666666666664
$ws.Columns.Item(5).ColumnWidth = 58.166666666666664
